$d = $word.ActiveDocument

# Namespace-qualified WordprocessingML fragments used to replace whole
# paragraphs in one shot via Range.InsertXML (keeps pPr + splits runs
# exactly the way we want, instead of letting same-format runs coalesce).
function Set-ParagraphXml($paraIndex, $innerXml) {
    $p = $d.Paragraphs.Item($paraIndex)
    $r = $p.Range
    $xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + $innerXml + '</w:p>'
    $null = $r.InsertXML($xml)
}

$plainPPr = '<w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/></w:pPr>'

# The six "Ingredients" list items: drop the w:ind left=720, drop any
# "(optional)" suffix, and split the text into 3 runs - a 16-space
# indent run, the item name, and a literal "\n" run.
$items = @(
    @{ Index = 2; Name = "Beets" },
    @{ Index = 3; Name = "Knife" },
    @{ Index = 4; Name = "Box grater" },
    @{ Index = 5; Name = "Food dehydrator" },
    @{ Index = 6; Name = "Cooling or baking rack" },
    @{ Index = 7; Name = "Cheesecloth" }
)

foreach ($item in $items) {
    $runs = '<w:r><w:t xml:space="preserve">                </w:t></w:r>' + `
            '<w:r><w:t>' + $item.Name + '</w:t></w:r>' + `
            '<w:r><w:t>\n</w:t></w:r>'
    $full = $plainPPr + $runs
    Set-ParagraphXml $item.Index $full
}

# The two instruction paragraphs whose trailing "\n" run gets merged
# back into the sentence run. (Range.Text = "..." only overwrites the
# first run in a multi-run range in this host, so go through InsertXML
# here too to fully collapse both runs into one.)
$firstLinePPr = '<w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:firstLine="720"/></w:pPr>'

$runs22 = '<w:r><w:t>If using the oven, dry beet strips for 4 hours and shredded beets for 2 hours.\n</w:t></w:r>'
$full22 = $firstLinePPr + $runs22
Set-ParagraphXml 22 $full22

$runs24 = '<w:r><w:t>Pack the dried and cooled beets into dry, airtight containers or plastic bags that seal.\n</w:t></w:r>'
$full24 = $firstLinePPr + $runs24
Set-ParagraphXml 24 $full24
